$d = $word.ActiveDocument

$replacements = @(
    @{ old = "981×7=6867"; new = "957×2=1914" },
    @{ old = "926×7=6482"; new = "660×6=3960" },
    @{ old = "472×4=1888"; new = "759×5=3795" },
    @{ old = "875×3=2625"; new = "556×3=1668" },
    @{ old = "476×7=3332"; new = "784×4=3136" },
    @{ old = "693×8=5544"; new = "861×5=4305" },
    @{ old = "338×3=1014"; new = "109×4=436" },
    @{ old = "986×7=6902"; new = "351×4=1404" },
    @{ old = "702×2=1404"; new = "725×3=2175" },
    @{ old = "448×4=1792"; new = "750×6=4500" },
    @{ old = "775×2=1550"; new = "756×3=2268" },
    @{ old = "177×2=354";  new = "924×4=3696" },
    @{ old = "285×9=2565"; new = "837×6=5022" },
    @{ old = "684×8=5472"; new = "322×8=2576" },
    @{ old = "802×5=4010"; new = "234×5=1170" },
    @{ old = "206×8=1648"; new = "313×5=1565" },
    @{ old = "761×2=1522"; new = "916×8=7328" },
    @{ old = "498×8=3984"; new = "280×2=560" },
    @{ old = "882×5=4410"; new = "558×2=1116" },
    @{ old = "427×7=2989"; new = "436×8=3488" },
    @{ old = "145×8=1160"; new = "300×9=2700" },
    @{ old = "679×6=4074"; new = "860×2=1720" },
    @{ old = "121×7=847";  new = "782×7=5474" },
    @{ old = "829×8=6632"; new = "910×5=4550" },
    @{ old = "540×7=3780"; new = "272×6=1632" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
